# Update the "placesToGo" sheet: replace "Miami" with "Kansas City" in A4
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("placesToGo")

$ws.Range("A4").Value = "Kansas City"

# Move the active selection to A5 (mirrors pressing Enter after editing A4)
$ws.Range("A5").Select()
